$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.473.59'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.650.26'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.76'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.84'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.85%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.111'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.390'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +8.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.73'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.87'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.127.96'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.343.87'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.637.54'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.24'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +8.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.67'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '354.30'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.98'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.75%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.73'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.08'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.74'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +11.69%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +9.23%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +6.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.27'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.54%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '543.07'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.07'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0₃0873'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +8.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.77'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.35'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.79%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.412'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.00'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.70'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '167.09'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.28'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.95'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +5.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0585'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.85'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.634'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.04'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +15.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0248'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0969'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.51'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.83%  '
